$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: Ativação (Activation date): 01/01/2020 -> 01/01/2022 ---
# The new value looks like a date, and a plain .Value assignment would make
# Excel auto-convert it to a date serial (changing cell type/format). Force
# the cells to Text first so it is stored as a literal string, matching the
# original authoring (a plain shared string, not a date), then restore the
# original "General" number format/style by pasting the format from the
# untouched neighboring row (which still carries the original style).
$ws.Range("B8").NumberFormat = "@"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("B8").Value = "01/01/2022"
$ws.Range("C8").Value = "01/01/2022"
$ws.Range("B9").Copy()
$ws.Range("B8").PasteSpecial(-4122)
$ws.Range("C9").Copy()
$ws.Range("C8").PasteSpecial(-4122)

# --- Row 16: Programa resumido (short program) ---
$ws.Range("B16").Value = "Estática de Partículas. Estática de Corpos Rígidos. Equilíbrio de Corpos Rígidos. Análise de Estruturas."
$ws.Range("C16").Value = "Estática de Partículas. Estática de Corpos Rígidos. Equilíbrio de Corpos Rígidos. Análise de Estruturas."

# --- Row 21: Método (evaluation method) ---
$ws.Range("B21").Value = "Os alunos serão avaliados continuamente a qual serão considerados: provas escritas, exercícios propostos e seminários"
$ws.Range("C21").Value = "Os alunos serão avaliados continuamente a qual serão considerados: provas escritas, exercícios propostos e seminários"

# --- Row 22: Critério (grading criteria) ---
$ws.Range("B22").Value = "Para compor a Nota Final (NF) serão consideradas as avaliações de Provas Escritas (P1 e P2) e Exercícios Propostos e Seminários (T) em que:NF = (P1 + P2 + T)/3.  Serão considerados aprovados os alunos que obtiverem: NF maior ou igual a 5,0. Serão considerados reprovados os alunos que obtiverem: NS menor que 3,0 Para os alunos em que NS é maior ou igual a 3,0 e menor que 5,0 será dada uma prova de recuperação (R)."
$ws.Range("C22").Value = "Para compor a Nota Final (NF) serão consideradas as avaliações de Provas Escritas (P1 e P2) e Exercícios Propostos e Seminários (T) em que:NF = (P1 + P2 + T)/3.  Serão considerados aprovados os alunos que obtiverem: NF maior ou igual a 5,0. Serão considerados reprovados os alunos que obtiverem: NS menor que 3,0 Para os alunos em que NS é maior ou igual a 3,0 e menor que 5,0 será dada uma prova de recuperação (R)."

# --- Row 23: Norma de recuperação (recovery rule) ---
$ws.Range("B23").Value = ": A prova de Recuperação (R) irá compor a nota final de recuperação (NR) da seguinte forma: NR = (R + NF)/2. Serão considerados aprovados os alunos que obtiverem NF maior ou igual a 5,0."
$ws.Range("C23").Value = ": A prova de Recuperação (R) irá compor a nota final de recuperação (NR) da seguinte forma: NR = (R + NF)/2. Serão considerados aprovados os alunos que obtiverem NF maior ou igual a 5,0."

# --- Row 24: Bibliografia (bibliography) ---
$ws.Range("B24").Value = "1. F.P. BEER, E.R. JOHNSTON, J.T. DeWOLF, D.MAZUREK. Estática e Mecânica dos Materiais. São Paulo: McGraw Hill, 2013, 728p.2. F.P. BEER, E.R. JOHNSTON, E. RUSSEL. Mecânica vetorial para engenheiros: Estática. São Paulo: McGraw Hill. 9a Ed., 2012, 626p. 3. HIBBELER, R.C. Mecânica para engenharia vol.1: estática. São Paulo: Pearson Prentice Hall, 12a Ed., 2011.4. MERIAM, J.L. KRAIGE, L.G. Mecânica para engenharia – Estática. Grupo GEN Editora LTC, 6a Ed., 2009, 364p. 5. RUIZ, C.C.de La P. Fundamentos de mecânica para engenharia – Estática. Grupo GEN Editora LTC, 2017, 306p."
$ws.Range("C24").Value = "1. F.P. BEER, E.R. JOHNSTON, J.T. DeWOLF, D.MAZUREK. Estática e Mecânica dos Materiais. São Paulo: McGraw Hill, 2013, 728p.2. F.P. BEER, E.R. JOHNSTON, E. RUSSEL. Mecânica vetorial para engenheiros: Estática. São Paulo: McGraw Hill. 9a Ed., 2012, 626p. 3. HIBBELER, R.C. Mecânica para engenharia vol.1: estática. São Paulo: Pearson Prentice Hall, 12a Ed., 2011.4. MERIAM, J.L. KRAIGE, L.G. Mecânica para engenharia – Estática. Grupo GEN Editora LTC, 6a Ed., 2009, 364p. 5. RUIZ, C.C.de La P. Fundamentos de mecânica para engenharia – Estática. Grupo GEN Editora LTC, 2017, 306p."
